$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (grandes-grupos) metadata re-curated: dimension -> measure
$ws.Range("D2").Value = "iaest-measure:grandes-grupos"
$ws.Range("D3").Value = "medida"
$ws.Range("D4").Value = "xsd:int"

# Column E (municipio-nombre) metadata re-curated
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Municipio"

# Column H (sexo) metadata re-curated: dimension -> measure
$ws.Range("H2").Value = "iaest-measure:sexo"
$ws.Range("H3").Value = "medida"
$ws.Range("H4").Value = "xsd:int"

# Row 5 (old mapping-file references) is no longer needed
$ws.Rows(5).Delete()
